# Re-shape the small wide-format "Year" table (one row per year, one
# column per metric) into a long/unpivoted layout (one row per metric,
# one column per year) with generic pandas-style "Unnamed: N" headers.
#
# Original layout (A1:K4):
#   Row1: Year | Flexible working hours | Full-time | Of which: female | ...
#   Row2: 2019 | 100 | 78.9 | 58.4 | ...
#   Row3: 2018 | 100 | 79.9 | 59.9 | ...
#   Row4: 2017 | 100 | 80.8 |      | ...
#
# New layout (A1:D12):
#   Row1 : Unnamed: 0 | Unnamed: 1 | Unnamed: 2 | Unnamed: 3
#   Row2 : (%)        | 2019       | 2018       | 2017
#   Row3 : Flexible working hours           | 100  | 100  | 100
#   Row4 : Full-time                        | 78.9 | 79.9 | 80.8
#   Row5 : Of which: female                 | 58.4 | 59.9 |
#   Row6 : Of which: male                   | 93.9 | 94.6 |
#   Row7 : Part-time                        | 21.1 | 20.1 | 19.2
#   Row8 : Of which: female                 | 41.6 | 40.1 |
#   Row9 : Of which: male                   | 6.1  | 5.4  |
#   Row10: Virtual offices                  | 13.2 | 12.8 | 12.9
#   Row11: Sabbatical                       | 0    | 0    | 0
#   Row12: Semi-retirement (Altersteilzeit) | 1.76 | 1.68 | 1.8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Snapshot the old wide-format values before anything is overwritten.
# ---------------------------------------------------------------------
$years = @(
    $ws.Range("A2").Value2,
    $ws.Range("A3").Value2,
    $ws.Range("A4").Value2
)

# Old column letters B..K, in order, correspond to these row labels in
# the new, unpivoted table.
$oldCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
$labels  = @(
    "Flexible working hours",
    "Full-time",
    "Of which: female",
    "Of which: male",
    "Part-time",
    "Of which: female",
    "Of which: male",
    "Virtual offices",
    "Sabbatical",
    "Semi-retirement (Altersteilzeit)"
)

$data = @{}
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $col = $oldCols[$i]
    $data[$i] = @(
        $ws.Range($col + "2").Value2,
        $ws.Range($col + "3").Value2,
        $ws.Range($col + "4").Value2
    )
}

# ---------------------------------------------------------------------
# 2. Wipe all old cell *contents* first (keeping formatting in place, so
#    A1's existing bold/bordered/centered style survives), then drop the
#    old E:K columns completely since the new layout only needs A:D.
# ---------------------------------------------------------------------
$ws.Range("A1:K4").ClearContents() | Out-Null
$ws.Range("E1:K4").Clear() | Out-Null

# ---------------------------------------------------------------------
# 3. Write the new text labels: generic header names on row 1, then the
#    metric name down column A for every data row.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"
$ws.Range("D1").Value = "Unnamed: 3"

$ws.Range("A2").Value = "(%)"
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item($i + 3, 1).Value = $labels[$i]
}

# ---------------------------------------------------------------------
# 4. Clone A1's pre-existing header style onto B1:D1 via copy/paste
#    special (keeps the style table untouched instead of growing it
#    with new, equivalent entries).
# ---------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# A2:A4 previously carried that same header-like style (bold/bordered/
# centered) because they held the year values; in the new layout only
# row 1 keeps it, so strip A2:A4 back to the default formatting.
$ws.Range("A2:A4").ClearFormats() | Out-Null

# ---------------------------------------------------------------------
# 5. Fill in the numeric values: row 2 gets the three years, rows 3-12
#    get the metric values per year (cells with no original value, e.g.
#    the "Of which" rows which only had 2019/2018 data, stay blank).
# ---------------------------------------------------------------------
$ws.Range("B2").Value = $years[0]
$ws.Range("C2").Value = $years[1]
$ws.Range("D2").Value = $years[2]

for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 3
    $vals = $data[$i]
    if ($null -ne $vals[0]) { $ws.Cells.Item($r, 2).Value = $vals[0] }
    if ($null -ne $vals[1]) { $ws.Cells.Item($r, 3).Value = $vals[1] }
    if ($null -ne $vals[2]) { $ws.Cells.Item($r, 4).Value = $vals[2] }
}
